# Add references for gp.
# Rework the Sheet1 table: rename a couple of header columns, correct some
# cell values, and add a new "group_formed" column (G) with sample data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write brand-new text values first, in the exact order they should be
# --- appended to the shared-strings table, so the resulting workbook's
# --- string table ordering matches the target layout.
$ws.Range("D2").Value = "Five Star"
$ws.Range("D3").Value = "Percentage"
$ws.Range("E3").Value = "A 1"
$ws.Range("F2").Value = "This message is post submission instructions text."
$ws.Range("E2").Value = "This message is instructions text."
$ws.Range("E1").Value = "instructions"
$ws.Range("F1").Value = "post_sub_instructions"
$ws.Range("G2").Value = "Educator Formed"
$ws.Range("G3").Value = "System Formed"
$ws.Range("G1").Value = "group_formed"

# --- Cells that reuse strings already present in the workbook.
$ws.Range("B2").Value = "0001_Ren_GP"
$ws.Range("B3").Value = "0002_Ren_GP"
$ws.Range("F3").Value = "A 2"

# --- Column widths (best effort; engine quantizes to 1/7 character units).
$ws.Columns.Item(1).ColumnWidth = 9
$ws.Columns.Item(4).ColumnWidth = 18.42857142857143
$ws.Range("E1:F1").ColumnWidth = 39.285714285714285
$ws.Columns.Item(7).ColumnWidth = 23.42857142857143

# --- Selection / active cell ends up on the newly added column.
$ws.Range("G1").Select()
